$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 4 corresponds to the f980c316... file.
# D4 = Correspond Handoff Datetime, G4 = Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-28 09:03:05"
$wsZhCn.Range("G4").Value = "2016-01-28 09:03:55"

# de-de sheet: row 4 corresponds to the f980c316... file.
# D4 = Correspond Handoff Datetime, G4 = Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-28 09:03:18"
$wsDeDe.Range("G4").Value = "2016-01-28 09:04:17"
